$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 0
$ws1.Range("F4").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("F10").Value = 0
$ws1.Range("F12").Value = 68

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 14
$ws2.Range("F3").Value = 20
$ws2.Range("F5").Value = 0
$ws2.Range("F6").Value = 0

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1162
$ws4.Range("F3").Value = 0
$ws4.Range("F4").Value = 345
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 20
$ws4.Range("F7").Value = 4963
$ws4.Range("F10").Value = 9193
$ws4.Range("F11").Value = 234
$ws4.Range("F13").Value = 79
